$wb = $excel.ActiveWorkbook

# --- 1. Status text update: "Ready for handoff" -> "In Translation" ---
# Overview sheet holds per-language status in columns E (zh-cn) and F (de-de), row 2.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn / de-de sheets hold the same status in column C, row 2.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Column width adjustments (report regenerated -> columns re-sized to new content) ---
# Overview: columns E and F (zh-cn / de-de status columns) shrink.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C (Status) shrinks the same way.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
